# Auto-generated edit script: updates "想去人数" (want-to-go count) values
# in worksheets "展览" (sheet1 / index 1) and "全部类型" (sheet4 / index 4)
# to match the data snapshot captured at commit 456a3b4 (gh-pages output).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (Exhibition) sheet updates ---
$wsExhibition.Range("F2").Value = 7701
$wsExhibition.Range("F3").Value = 3618
$wsExhibition.Range("F5").Value = 76
$wsExhibition.Range("F6").Value = 105
$wsExhibition.Range("F7").Value = 91
$wsExhibition.Range("F8").Value = 127
$wsExhibition.Range("F10").Value = 535
$wsExhibition.Range("F15").Value = 21
$wsExhibition.Range("F17").Value = 367
$wsExhibition.Range("F18").Value = 4307
$wsExhibition.Range("F19").Value = 4307
$wsExhibition.Range("F20").Value = 123
$wsExhibition.Range("F22").Value = 1048
$wsExhibition.Range("F24").Value = 2684
$wsExhibition.Range("F27").Value = 3148
$wsExhibition.Range("F28").Value = 2452
$wsExhibition.Range("F29").Value = 80
$wsExhibition.Range("F32").Value = 107
$wsExhibition.Range("F33").Value = 138
$wsExhibition.Range("F34").Value = 145
$wsExhibition.Range("F37").Value = 116
$wsExhibition.Range("F38").Value = 4612
$wsExhibition.Range("F39").Value = 573
$wsExhibition.Range("F40").Value = 348
$wsExhibition.Range("F43").Value = 903
$wsExhibition.Range("F46").Value = 1739
$wsExhibition.Range("F47").Value = 273
$wsExhibition.Range("F49").Value = 633
$wsExhibition.Range("F50").Value = 751

# --- 全部类型 (All types) sheet updates ---
$wsAllTypes.Range("F3").Value = 7701
$wsAllTypes.Range("F4").Value = 3618
$wsAllTypes.Range("F6").Value = 76
$wsAllTypes.Range("F7").Value = 105
$wsAllTypes.Range("F8").Value = 91
$wsAllTypes.Range("F9").Value = 127
$wsAllTypes.Range("F12").Value = 535
$wsAllTypes.Range("F16").Value = 21
$wsAllTypes.Range("F17").Value = 367
$wsAllTypes.Range("F18").Value = 4307
$wsAllTypes.Range("F19").Value = 4307
$wsAllTypes.Range("F24").Value = 1048
$wsAllTypes.Range("F26").Value = 2684
$wsAllTypes.Range("F29").Value = 3148
$wsAllTypes.Range("F30").Value = 2452
$wsAllTypes.Range("F31").Value = 80
$wsAllTypes.Range("F33").Value = 138
$wsAllTypes.Range("F34").Value = 145
$wsAllTypes.Range("F37").Value = 116
$wsAllTypes.Range("F39").Value = 4612
$wsAllTypes.Range("F41").Value = 573
$wsAllTypes.Range("F42").Value = 348
$wsAllTypes.Range("F45").Value = 903
$wsAllTypes.Range("F47").Value = 1739
$wsAllTypes.Range("F48").Value = 273
$wsAllTypes.Range("F49").Value = 633
$wsAllTypes.Range("F50").Value = 751

Write-Host "Done updating 想去人数 (want-to-go counts)."
